# CS133JS Lab07 Instructions - groupA : "Minor updates and clarifications"
#
# This script reproduces (via Word COM-interop calls) the following
# source edits to the three paragraphs under the "Upload the following..."
# heading at the end of the document:
#
#   1. "Upload the following to the Lab Production Version assignment:"
#        -> "Upload the following 6 files to the Lab Production Version assignment:"
#      (split into extra runs; the list items' leading bookmark "_GoBack"
#       moves here, to just before the first run of this paragraph)
#
#   2. "A zip file containing the four files for part 2."
#        -> "The four files for part 2."
#      (the old "_GoBack" bookmarkStart/bookmarkEnd pair that used to sit
#       at the top of this paragraph is removed from here)
#
#   3. "A code review of your own code."
#        -> "A code review of your code with the \u201cProd\u201d column filled in by you."
#      (the "_GoBack" bookmarkEnd now closes at the very end of this
#       paragraph's content)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: force a just-inserted run of text to become its own <w:r>
# element (rather than being silently re-merged into a neighbouring run
# that happens to carry identical formatting) by toggling a character
# property on/off. Needs to run *after* all sibling insertions in the
# same paragraph/region are done, and right-to-left (rightmost /
# most-recently-inserted text first) so later toggles don't cause the
# engine to re-normalize (merge) runs it already split.
# ---------------------------------------------------------------------
function Split-Run($start, $end) {
    $rng = $d.Range($start, $end)
    $rng.Bold = 1
    $rng.Bold = 0
}

# =======================================================================
# Paragraph: "Upload the following to the Lab Production Version assignment:"
# =======================================================================

$pUpload = $d.Paragraphs.Item(42)
$uploadStart = $pUpload.Range.Start

# New bookmark start goes immediately before this paragraph's first run.
$bmStartRange = $d.Range($uploadStart, $uploadStart)

$fr = $d.Range($pUpload.Range.Start, $pUpload.Range.End)
$fr.Find.Execute("Upload the following to the ") | Out-Null
$fr.Text = "Upload the following "

$fr.Collapse(0)
$fr.InsertAfter("6 files ")
$run6filesStart = $fr.Start
$run6filesEnd = $fr.End

$fr.Collapse(0)
$fr.InsertAfter("to the ")
$runToTheStart = $fr.Start
$runToTheEnd = $fr.End

Split-Run $runToTheStart $runToTheEnd
Split-Run $run6filesStart $run6filesEnd

# =======================================================================
# Paragraph: "A zip file containing the four files for part 2."
# =======================================================================

$pZip = $d.Paragraphs.Item(44)
$frZip = $d.Range($pZip.Range.Start, $pZip.Range.End)
$frZip.Find.Execute("A zip file containing the four files") | Out-Null
$frZip.Text = "The"

$frZip.Collapse(0)
$frZip.InsertAfter(" four files")
$runFourFilesStart = $frZip.Start
$runFourFilesEnd = $frZip.End

Split-Run $runFourFilesStart $runFourFilesEnd

# =======================================================================
# Paragraph: "A code review of your own code."
# =======================================================================

$pCode = $d.Paragraphs.Item(45)
$runAStart = $pCode.Range.Start

$frCode = $d.Range($pCode.Range.Start, $pCode.Range.End)
$frCode.Find.Execute(" code review ") | Out-Null
$runAEnd = $frCode.Start
$frCode.Text = " code review"
$runCodeReviewStart = $frCode.Start
$runCodeReviewEnd = $frCode.End

$frCode.Collapse(0)
$frCode.InsertAfter(" of your code ")
$runOfYourCodeStart = $frCode.Start
$runOfYourCodeEnd = $frCode.End

$frCode.Collapse(0)
$frCode.InsertAfter("with the ")
$runWithTheStart = $frCode.Start
$runWithTheEnd = $frCode.End

Split-Run $runWithTheStart $runWithTheEnd
Split-Run $runOfYourCodeStart $runOfYourCodeEnd
Split-Run $runCodeReviewStart $runCodeReviewEnd
Split-Run $runAStart $runAEnd

# "of your own code." (underlined run) -> "\u201cProd\u201d column filled in by you" + "."
# (the trailing "." loses the underline formatting)
$pCode2 = $d.Paragraphs.Item(45)
$frProd = $d.Range($pCode2.Range.Start, $pCode2.Range.End)
$frProd.Find.Execute("of your own code.") | Out-Null
$frProd.Text = "“Prod” column filled in by you"

# Merge a "." onto the end of the paragraph (inherits the underlined run's
# formatting because it is inserted adjacent to existing text, not at the
# bare end-of-story position), then strip the underline back off of just
# that trailing character - this naturally forces it into its own run.
$pCode3 = $d.Paragraphs.Item(45)
$paraEnd = $pCode3.Range.End
$markRange = $d.Range($paraEnd - 1, $paraEnd)
$markRange.InsertBefore(".")

$pCode4 = $d.Paragraphs.Item(45)
$paraEnd2 = $pCode4.Range.End
$periodRange = $d.Range($paraEnd2 - 2, $paraEnd2 - 1)
$periodRange.Font.Underline = 0

# =======================================================================
# Bookmark "_GoBack": remove its old location (top of the "zip file"
# paragraph) and re-add it spanning from the start of the "Upload the
# following" paragraph through to the end of the "code review" paragraph.
# =======================================================================

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$pCodeFinal = $d.Paragraphs.Item(45)
$bmEndPos = $pCodeFinal.Range.End - 1
$bmRange = $d.Range($uploadStart, $bmEndPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
